# ZPD Latvija Arpus Civilizacijas Kopsavilkums: add Bx/By/Bz STD + error columns
# (L:N = error, O:Q = STD) for rows 1-11 on the "Results" sheet, and reposition
# the two scatter charts that were moved further down/right on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header labels (L1:N1 = errors, O1:Q1 = STD) -----------------------
# Written in STD-then-error order so the shared-string table indices line up
# with the target workbook (91=Bx STD,92=By STD,93=Bz STD,94=Bx error,
# 95=By error,96=Bz error).
$ws.Cells.Item(1, 15).Value = "Bx STD"
$ws.Cells.Item(1, 16).Value = "By STD"
$ws.Cells.Item(1, 17).Value = "Bz STD"
$ws.Cells.Item(1, 12).Value = "Bx error"
$ws.Cells.Item(1, 13).Value = "By error"
$ws.Cells.Item(1, 14).Value = "Bz error"

# --- New data values (rows 2-11) --------------------------------------------
$ws.Cells.Item(2, 12).Value = 1.2303471472024561
$ws.Cells.Item(2, 13).Value = 0.96249760852452781
$ws.Cells.Item(2, 14).Value = 0.75302631850449531
$ws.Cells.Item(2, 15).Value = 143.73488091859826
$ws.Cells.Item(2, 16).Value = 112.44345098883198
$ws.Cells.Item(2, 17).Value = 87.972039814063635
$ws.Cells.Item(3, 12).Value = 0.63228473095822002
$ws.Cells.Item(3, 13).Value = 0.74957980295904392
$ws.Cells.Item(3, 14).Value = 0.80492609956858752
$ws.Cells.Item(3, 15).Value = 67.138397253790586
$ws.Cells.Item(3, 16).Value = 79.593234061202992
$ws.Cells.Item(3, 17).Value = 85.470114312077172
$ws.Cells.Item(4, 12).Value = 0.62693378268198097
$ws.Cells.Item(4, 13).Value = 0.27240676173106815
$ws.Cells.Item(4, 14).Value = 0.71563687811763221
$ws.Cells.Item(4, 15).Value = 75.900426553744467
$ws.Cells.Item(4, 16).Value = 32.979223616029522
$ws.Cells.Item(4, 17).Value = 86.63936416754126
$ws.Cells.Item(5, 12).Value = 1.1352250317262729
$ws.Cells.Item(5, 13).Value = 0.97615341563397662
$ws.Cells.Item(5, 14).Value = 0.90333535191682945
$ws.Cells.Item(5, 15).Value = 140.12098707882188
$ws.Cells.Item(5, 16).Value = 120.48675488682916
$ws.Cells.Item(5, 17).Value = 111.49881092852904
$ws.Cells.Item(6, 12).Value = 0.6156152455628866
$ws.Cells.Item(6, 13).Value = 1.3526776234906215
$ws.Cells.Item(6, 14).Value = 0.734797209996773
$ws.Cells.Item(6, 15).Value = 67.279733705786526
$ws.Cells.Item(6, 16).Value = 147.8322555430095
$ws.Cells.Item(6, 17).Value = 80.304964785489034
$ws.Cells.Item(7, 12).Value = 0.62261802504941355
$ws.Cells.Item(7, 13).Value = 0.33656647938121664
$ws.Cells.Item(7, 14).Value = 0.69658342722684652
$ws.Cells.Item(7, 15).Value = 67.993765333437523
$ws.Cells.Item(7, 16).Value = 36.755155323894584
$ws.Cells.Item(7, 17).Value = 76.071247828497903
$ws.Cells.Item(8, 12).Value = 0.6393485929415561
$ws.Cells.Item(8, 13).Value = 0.31898794729853674
$ws.Cells.Item(8, 14).Value = 0.703024193225951
$ws.Cells.Item(8, 15).Value = 66.933416103826417
$ws.Cells.Item(8, 16).Value = 33.394854144287038
$ws.Cells.Item(8, 17).Value = 73.599615883648141
$ws.Cells.Item(9, 12).Value = 0.63132744384288686
$ws.Cells.Item(9, 13).Value = 1.5444318399303869
$ws.Cells.Item(9, 14).Value = 0.91921464451443335
$ws.Cells.Item(9, 15).Value = 66.929642330013877
$ws.Cells.Item(9, 16).Value = 163.73162874153533
$ws.Cells.Item(9, 17).Value = 97.449759204785209
$ws.Cells.Item(10, 12).Value = 0.63764837653464124
$ws.Cells.Item(10, 13).Value = 1.7652044788642125
$ws.Cells.Item(10, 14).Value = 0.99849452898145585
$ws.Cells.Item(10, 15).Value = 68.525695007575791
$ws.Cells.Item(10, 16).Value = 189.69994780200676
$ws.Cells.Item(10, 17).Value = 107.30448641861966
$ws.Cells.Item(11, 12).Value = 0.68623044829777413
$ws.Cells.Item(11, 13).Value = 0.78521469712324965
$ws.Cells.Item(11, 14).Value = 0.74684855950169626
$ws.Cells.Item(11, 15).Value = 69.583902808435894
$ws.Cells.Item(11, 16).Value = 79.620925162840592
$ws.Cells.Item(11, 17).Value = 75.730591240736914

# --- Reposition the two scatter charts ---------------------------------------
# (they got dragged further down/right and slightly resized; positions below
# are derived from the target col/row/offset anchors converted to points
# using this sheet's actual column widths / 14.5pt default row height)
$charts = $ws.ChartObjects()

$chart1 = $charts.Item(1)
$chart1.Left = 779.0811629552
$chart1.Top = 224.2960629921
$chart1.Width = 528.2168700787
$chart1.Height = 303.4512598425

$chart2 = $charts.Item(2)
$chart2.Left = 716.9724818529
$chart2.Top = 547.3681102362
$chart2.Width = 528.0046653543
$chart2.Height = 269.9655905512

# --- Selection / scroll position ---------------------------------------------
$ws.Range("H1").Select()
